# checkpoint grupo 13; feito g13.5a
#
# Swap the Rio Grande do Norte / Distrito Federal rows (4 & 5) - both share
# the same 8.8 value, so they now tie for 3rd place instead of 3rd/4th -
# and swap the Nordeste / Brasil rows (8 & 9). Also add a thin box border +
# top vertical alignment to the header row style, and reset the page
# margins to Excel's stock defaults.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 4 & 5: Distrito Federal and Rio Grande do Norte trade places ---
# (only the region name and ranking text actually move; the value is the
# same 8.8 for both rows, and the variable/quarter columns are untouched)
$ws.Range("A4").Value = "Distrito Federal"
$ws.Range("E4").Value = "3º"

$ws.Range("A5").Value = "Rio Grande do Norte"
$ws.Range("E5").Value = "3º"

# --- Rows 8 & 9: Brasil / Nordeste swap places (region name + value only) ---
$ws.Range("A8").Value = "Brasil"
$ws.Range("D8").Value = 6.4

$ws.Range("A9").Value = "Nordeste"
$ws.Range("D9").Value = 8.699999999999999

# --- Header row (A1:E1): add a thin box border and top vertical alignment ---
$headerRng = $ws.Range("A1:E1")
$headerRng.Borders.LineStyle = 1
$headerRng.VerticalAlignment = -4160

# --- Page margins: back to Excel's stock defaults (points) ---
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36
